$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'45.310.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'318.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.49%  "
$ws.Range("D6").Value = "'102.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.15%  "
$ws.Range("D7").Value = "'0.516"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +5.27%  "
$ws.Range("D10").Value = "'35.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("E12").Value = "  -1.85%  "
$ws.Range("D13").Value = "'18.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.25%  "
$ws.Range("E14").Value = "  +2.16%  "
$ws.Range("D15").Value = "'2.808.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "'2.449.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("D17").Value = "'0.846"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").Value = "'45.224.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.32%  "
$ws.Range("D19").Value = "'12.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.95%  "
$ws.Range("D20").Value = "'6.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E21").Value = "  +1.61%  "
$ws.Range("D22").Value = "'68.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("D23").Value = "'245.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'2.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.27%  "
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'25.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.12%  "
$ws.Range("E28").Value = "  -2.70%  "
$ws.Range("D29").Value = "'9.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("D30").Value = "'49.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.65%  "
$ws.Range("D31").Value = "'33.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("D32").Value = "'0.126"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.19%  "
$ws.Range("D33").Value = "'20.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.99%  "
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").Value = "'0.0764"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("D40").Value = "'126.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.52%  "
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("E42").Value = "  -4.15%  "
$ws.Range("D43").Value = "'20.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.32%  "
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("D45").Value = "'1.937.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("E46").Value = "  -2.73%  "
$ws.Range("E47").Value = "  +2.02%  "
$ws.Range("E48").Value = "  +8.51%  "
$ws.Range("D49").Value = "'9.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.45%  "
$ws.Range("D50").Value = "'76.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.25%  "
$ws.Range("D51").Value = "'4.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.44%  "
